# Fruta / hortaliza, semanal
# Updates the "Fruta, Vega Modelo de Temuco - Kiwi" dataset:
#  - Rows 291-293 get new weekly price entries (date 2021-09-20 / serial 44448)
#    replacing their previous data, while the displaced original rows are
#    re-appended further down the table (now rows 298-301).
#  - Rows 294-297 are brand-new entries for the same reporting date.
#  - The sheet dimension grows from A1:T294 to A1:T301.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column letter -> index map (A=1 ... T=20)
$colIndex = @{
    A=1; B=2; C=3; D=4; E=5; F=6; G=7; H=8; I=9; J=10;
    K=11; L=12; M=13; N=14; O=15; P=16; Q=17; R=18; S=19; T=20
}

# Use the existing date-formatted cell (column D) as the template so new/moved
# date cells keep the same number format as the rest of the column.
$dateFormat = $ws.Range("D2").NumberFormat

# Final target state (after the edit) for every row from 291 through 301.
$rows = @(
    @{ Row=291; A=10; B="Vega Modelo de Temuco"; C="La Araucanía"; D=44448; E=9; F="Fruta"; G=100101; H="Berries"; I=100101007; J="Kiwi"; K="Hayward"; L="Especial"; M=130; N=22000; O=25000; P=23385; Q="`$/bandeja 18 kilos"; R="Región de O'Higgins"; S=1299; T=18 },
    @{ Row=292; A=10; B="Vega Modelo de Temuco"; C="La Araucanía"; D=44448; E=9; F="Fruta"; G=100101; H="Berries"; I=100101007; J="Kiwi"; K="Hayward"; L="Especial"; M=5; N=550000; O=550000; P=550000; Q="`$/bins (450 kilos)"; R="Región de O'Higgins"; S=1222; T=450 },
    @{ Row=293; A=10; B="Vega Modelo de Temuco"; C="La Araucanía"; D=44448; E=9; F="Fruta"; G=100101; H="Berries"; I=100101007; J="Kiwi"; K="Hayward"; L="Primera"; M=80; N=18000; O=18000; P=18000; Q="`$/bandeja 18 kilos"; R="Región de O'Higgins"; S=1000; T=18 },
    @{ Row=294; A=10; B="Vega Modelo de Temuco"; C="La Araucanía"; D=44448; E=9; F="Fruta"; G=100101; H="Berries"; I=100101007; J="Kiwi"; K="Hayward"; L="Primera"; M=8; N=400000; O=400000; P=400000; Q="`$/bins (450 kilos)"; R="Región de O'Higgins"; S=889; T=450 },
    @{ Row=295; A=10; B="Vega Modelo de Temuco"; C="La Araucanía"; D=44448; E=9; F="Fruta"; G=100101; H="Berries"; I=100101007; J="Kiwi"; K="Hayward"; L="Segunda"; M=140; N=13000; O=15000; P=14429; Q="`$/bandeja 18 kilos"; R="Región de O'Higgins"; S=802; T=18 },
    @{ Row=296; A=10; B="Vega Modelo de Temuco"; C="La Araucanía"; D=44448; E=9; F="Fruta"; G=100101; H="Berries"; I=100101007; J="Kiwi"; K="Hayward"; L="Segunda"; M=5; N=300000; O=300000; P=300000; Q="`$/bins (450 kilos)"; R="Región de O'Higgins"; S=667; T=450 },
    @{ Row=297; A=10; B="Vega Modelo de Temuco"; C="La Araucanía"; D=44448; E=9; F="Fruta"; G=100101; H="Berries"; I=100101007; J="Kiwi"; K="Hayward"; L="Tercera"; M=50; N=5000; O=5000; P=5000; Q="`$/caja 10 kilos"; R="Región de O'Higgins"; S=500; T=10 },
    @{ Row=298; A=10; B="Vega Modelo de Temuco"; C="La Araucanía"; D=44238; E=9; F="Fruta"; G=100101; H="Berries"; I=100101007; J="Kiwi"; K="Hayward"; L="Primera"; M=65; N=15000; O=15000; P=15000; Q="`$/bandeja 10 kilos"; R="Región de O'Higgins"; S=1500; T=10 },
    @{ Row=299; A=10; B="Vega Modelo de Temuco"; C="La Araucanía"; D=44399; E=9; F="Fruta"; G=100101; H="Berries"; I=100101007; J="Kiwi"; K="Hayward"; L="Primera"; M=235; N=9000; O=10000; P=9340; Q="`$/bandeja 10 kilos"; R="Región de O'Higgins"; S=934; T=10 },
    @{ Row=300; A=10; B="Vega Modelo de Temuco"; C="La Araucanía"; D=44399; E=9; F="Fruta"; G=100101; H="Berries"; I=100101007; J="Kiwi"; K="Hayward"; L="Primera"; M=250; N=20000; O=20000; P=20000; Q="`$/bandeja 18 kilos"; R="Región de O'Higgins"; S=1111; T=18 },
    @{ Row=301; A=10; B="Vega Modelo de Temuco"; C="La Araucanía"; D=44400; E=9; F="Fruta"; G=100101; H="Berries"; I=100101007; J="Kiwi"; K="Hayward"; L="Primera"; M=120; N=15000; O=16000; P=15458; Q="`$/bandeja 18 kilos"; R="Región de O'Higgins"; S=859; T=18 },
)

foreach ($r in $rows) {
    $rowNum = $r.Row
    foreach ($col in $colIndex.Keys) {
        $cell = $ws.Cells.Item($rowNum, $colIndex[$col])
        $cell.Value = $r[$col]
        if ($col -eq "D") {
            $cell.NumberFormat = $dateFormat
        }
    }
}

Write-Host "Updated rows 291-301 (now spanning A1:T301)."
